$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns K (Impressora HP Laser) and L (Impressora Multifuncional) ---
# Order matters for shared-string indices: Multifuncional must land before HP Laser.
$ws.Range("L1").Value = "Impressora Multifuncional"
$ws.Range("K1").Value = "Impressora HP Laser"
$ws.Range("K1").Style = $ws.Range("J1").Style

$ws.Range("L1").Style = $ws.Range("J1").Style

# --- Row 2 (Quantidade) ---
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1

# --- Row 3 (Valor) ---
$ws.Range("K3").Value = 1124.0999999999999
$ws.Range("L3").Value = 1093.4100000000001

# --- Row 4 (Total) ---
$ws.Range("K4").Value = 1124.0999999999999
$ws.Range("L4").Value = 1093.4100000000001

# --- Row 6 (Link row) with hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("L6"), "https://www.kalunga.com.br/prod/impressora-multifuncional-tanque-de-tinta-smart-tank-517-1tj10a-color-wi-fi-conexao-usb-bivolt-hp-cx-1-un/213106?pcID=84?utm_source=paid_search&utm_content=kalunga_hp&utm_medium=pla&utm_content=kalunga_hp&utm_source=google&utm_medium=cpc&utm_campaign=Kalunga_BR_MIX_Print_ConsHW_CISS_Local_OPEX_Google_All_Smart-PLA_ROAS&targetid=pla-1678689595881&utm_term=&matchtype=&adid=515801921680&addisttype=u&gclid=CjwKCAjwvNaYBhA3EiwACgndgq9Ik0-HxRHYHbcndhnhRuoxca4QhfXOe2wySSmvrnSrfCbdS9qelBoCwVwQAvD_BwE&gclsrc=aw.ds", "", "", "https://www.kalunga.com.br/prod/impressora-multifuncional-tanque-de-tinta-smart-tank-517-1tj10a-color-wi-fi-conexao-usb-bivolt-hp-cx-1-un/213106?pcID=84?utm_source=paid_search&utm_content=kalunga_hp&utm_medium=pla&utm_content=kalunga_hp&utm_source=google&utm_medium=cpc&utm_campaign=Kalunga_BR_MIX_Print_ConsHW_CISS_Local_OPEX_Google_All_Smart-PLA_ROAS&targetid=pla-1678689595881&utm_term=&matchtype=&adid=515801921680&addisttype=u&gclid=CjwKCAjwvNaYBhA3EiwACgndgq9Ik0-HxRHYHbcndhnhRuoxca4QhfXOe2wySSmvrnSrfCbdS9qelBoCwVwQAvD_BwE&gclsrc=aw.ds")
$ws.Hyperlinks.Add($ws.Range("K6"), "https://www.amazon.com.br/Impressora-Laser-Monocrom%C3%A1tica-HP-Branca/dp/B07S61ZJCS/ref=asc_df_B07S61ZJCS/?tag=googleshopp00-20&linkCode=df0&hvadid=379699119574&hvpos=&hvnetw=g&hvrand=2465213770169722673&hvpone=&hvptwo=&hvqmt=&hvdev=c&hvdvcmdl=&hvlocint=&hvlocphy=1001773&hvtargid=pla-779928478246&psc=1", "", "", "https://www.amazon.com.br/Impressora-Laser-Monocrom%C3%A1tica-HP-Branca/dp/B07S61ZJCS/ref=asc_df_B07S61ZJCS/?tag=googleshopp00-20&linkCode=df0&hvadid=379699119574&hvpos=&hvnetw=g&hvrand=2465213770169722673&hvpone=&hvptwo=&hvqmt=&hvdev=c&hvdvcmdl=&hvlocint=&hvlocphy=1001773&hvtargid=pla-779928478246&psc=1")

$ws.Range("K6").Style = $ws.Range("B6").Style
$ws.Range("L6").Style = $ws.Range("B6").Style

# --- Column widths ---
$ws.Columns("K").ColumnWidth = 28.28515625
$ws.Columns("L").ColumnWidth = 33.85546875

# --- Sheet view ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K7").Select()
